$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cell values from the crypto symbol list refresh.
# A leading apostrophe forces Excel to store the value as literal text
# (matching the workbook's original inlineStr string cells) instead of
# auto-converting numeric-looking / percentage-looking text into numbers.

$ws.Range('D2').Value = '''307.76'
$ws.Range('E2').Value = '''0.87%'
$ws.Range('E3').Value = '''1.19%'
$ws.Range('D4').Value = '''5.059'
$ws.Range('E4').Value = '''1.71%'
$ws.Range('D5').Value = '''0.08120'
$ws.Range('D6').Value = '''2.008'
$ws.Range('E6').Value = '''5.69%'
$ws.Range('B7').Value = '''GateToken'
$ws.Range('C7').Value = '''https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D7').Value = '''4.165'
$ws.Range('E7').Value = '''0.47%'
$ws.Range('B8').Value = '''KuCoinToken'
$ws.Range('C8').Value = '''https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range('D8').Value = '''7.859'
$ws.Range('E8').Value = '''-0.24%'
$ws.Range('B9').Value = '''MXToken'
$ws.Range('C9').Value = '''https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D9').Value = '''0.9278'
$ws.Range('E9').Value = '''-0.21%'
$ws.Range('B10').Value = '''LiechtensteinCryptoassetsExchange'
$ws.Range('C10').Value = '''https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D10').Value = '''0.1486'
$ws.Range('E10').Value = '''20.84%'
$ws.Range('B11').Value = '''WazirX'
$ws.Range('C11').Value = '''https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').Value = '''0.1929'
$ws.Range('E11').Value = '''1.30%'
$ws.Range('B12').Value = '''MandalaExchangeToken'
$ws.Range('C12').Value = '''https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').Value = '''0.09088'
$ws.Range('E12').Value = '''-1.45%'
$ws.Range('B13').Value = '''BitrueCoin'
$ws.Range('C13').Value = '''https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').Value = '''0.03514'
$ws.Range('E13').Value = '''0.22%'
$ws.Range('B14').Value = '''BitMartToken'
$ws.Range('C14').Value = '''https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').Value = '''0.09888'
$ws.Range('E14').Value = '''-0.24%'
$ws.Range('B15').Value = '''BitForexToken'
$ws.Range('C15').Value = '''https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').Value = '''0.001421'
$ws.Range('E15').Value = '''-0.30%'
$ws.Range('B16').Value = '''TigerCash'
$ws.Range('C16').Value = '''https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').Value = '''0.006101'
$ws.Range('E16').Value = '''-3.70%'
$ws.Range('B17').Value = '''LEO'
$ws.Range('C17').Value = '''https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').Value = '''3.841'
$ws.Range('E17').Value = '''6.43%'
$ws.Range('D19').Value = '''0.3457'
$ws.Range('D20').Value = '''0.1295'
$ws.Range('E20').Value = '''0.03%'
$ws.Range('D21').Value = '''4.841'
$ws.Range('E21').Value = '''-6.89%'
$ws.Range('E22').Value = '''-7.35%'
$ws.Range('D23').Value = '''0.04386'
$ws.Range('E23').Value = '''-0.63%'
$ws.Range('E24').Value = '''0.13%'
$ws.Range('D25').Value = '''0.004183'
$ws.Range('E25').Value = '''-11.26%'
$ws.Range('E27').Value = '''0.06%'
$ws.Range('D39').Value = '''0.02049'
$ws.Range('E39').Value = '''5.02%'
$ws.Range('D40').Value = '''0.05121'
$ws.Range('E40').Value = '''-1.51%'
$ws.Range('D41').Value = '''0.007488'
$ws.Range('E41').Value = '''-0.77%'
$ws.Range('D42').Value = '''0.01000'
$ws.Range('E42').Value = '''-1.57%'
$ws.Range('E43').Value = '''0.15%'
$ws.Range('D44').Value = '''0.002123'
$ws.Range('E44').Value = '''1.01%'
$ws.Range('D45').Value = '''0.009867'
$ws.Range('E45').Value = '''-7.90%'
$ws.Range('D46').Value = '''0.00006309'
$ws.Range('E46').Value = '''-0.57%'
$ws.Range('E47').Value = '''0.12%'
$ws.Range('D48').Value = '''63.84'
$ws.Range('E48').Value = '''0.42%'
$ws.Range('D49').Value = '''0.001602'
$ws.Range('E49').Value = '''-3.43%'
$ws.Range('E50').Value = '''0.12%'
$ws.Range('E51').Value = '''0.12%'
